$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column widths for the masker (F) / Soundfile (G) columns -----------
$ws.Columns.Item(6).ColumnWidth = 15
$ws.Columns.Item(7).ColumnWidth = 103.16666666666667

# --- new block of "perfectdata" trial rows (256-303) --------------------
# Row 256 seeds the counters, rows 257+ increment the cell above by 1
# (column B = Trial counter, column D = Block counter), mirroring a
# drag-filled series in column A of the literal "perfectdata" label.
$startRow = 256
$endRow = 303

$ws.Range("A$startRow").Value = "perfectdata"
$ws.Range("B$startRow").Value = 1
$ws.Range("D$startRow").Value = 1

for ($r = $startRow + 1; $r -le $endRow; $r++) {
    $prev = $r - 1
    $ws.Range("A$r").Value = "perfectdata"
    $ws.Range("B$r").Formula = "=B$prev+1"
    $ws.Range("D$r").Formula = "=D$prev+1"
}

# --- view state: select G170 (mirrors the author's final cursor spot) ---
$ws.Range("G170").Select() | Out-Null
